$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Red")

$ws.Range("C2").Value = "06.01.00"
$ws.Range("C3").Value = "06.10.00"
$ws.Range("C4").Value = "06.20.00"
$ws.Range("C5").Value = "06.30.00"
$ws.Range("C6").Value = "06.40.00"
$ws.Range("C7").Value = "06.50.00"

$ws.Range("D7").Select()
